# Swap the player-stat data between row 18 (John Stones) and row 19 (Marc Guéhi),
# keeping League/Team (A,B) and type/goalsPrevented (DK,DL) columns untouched,
# since those values are identical for both rows anyway.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row18Range = $ws.Range("C18:DJ18")
$row19Range = $ws.Range("C19:DJ19")

$row18Values = $row18Range.Value()
$row19Values = $row19Range.Value()

$row18Range.Value = $row19Values
$row19Range.Value = $row18Values
